# Generate Report for handoff
# Updates the status / latest-handoff info for "b.md.md" now that it is
# ready to be handed off again (new xlf files generated).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for b.md.md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row for b.md.md (row 3) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-02-15 03:40:55"
$wsZh.Hyperlinks.Item(2).TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"

# --- de-de sheet: row for b.md.md (row 3) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$wsDe.Range("D3").Value = "2016-02-15 03:41:09"
$wsDe.Hyperlinks.Item(2).TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
